$d = $word.ActiveDocument

# Locate the existing "Braubbuddy" project entry (DefinitionTerm paragraph
# containing the Braubbuddy hyperlink) -- the new SparkCC entry is inserted
# immediately before it.
$findRange = $d.Content
$findRange.Find.Execute("Braubbuddy")
$braubbuddyPara = $findRange.Paragraphs(1)

# Create two new empty paragraphs right before the Braubbuddy term paragraph;
# they inherit the DefinitionTerm style from the paragraph we split.
$insertionPoint = $braubbuddyPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()

# Re-locate the Braubbuddy paragraph (collection was shifted by the inserts)
# and walk backwards to the two fresh paragraphs.
$findRange2 = $d.Content
$findRange2.Find.Execute("Braubbuddy")
$braubbuddyPara2 = $findRange2.Paragraphs(1)
$descriptionPara = $braubbuddyPara2.Previous()
$termPara = $descriptionPara.Previous()

# First new paragraph: "SparkCC" hyperlink, styled like the other project
# links (DefinitionTerm paragraph style, Link character style).
$termPara.Style = "DefinitionTerm"
$termRange = $termPara.Range
$termRange.Collapse(1)
$sparkHyperlink = $d.Hyperlinks.Add($termRange, "http://sparkcc.org", "", "", "SparkCC")
$sparkHyperlink.Range.Style = "Link"

# Second new paragraph: plain description text.
$descriptionPara.Style = "Definition"
$descriptionRange = $descriptionPara.Range
$descriptionRange.Collapse(1)
$descriptionRange.InsertAfter("President and co-founder of SparkCC: a makerspace/hackerspace located on the Central Coast of New South Wales, Australia.")
